# Refresh the cryptos snapshot data: updated prices/volumes pulled from
# coinranking.com, plus a handful of rows that changed rank position
# (Chainlink/Polkadot and Frax/VeChain/MXToken swapped places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.039.47'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '1.922.42'
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''325.76'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('D8').Value = '''0.3817'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').Value = '''0.07750'
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('D10').Value = '''0.9781'
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('D11').Value = '''22.62'
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('D12').Value = '1.904.87'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '''6.965'
$ws.Range('E13').Value = '  +0.40%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''5.687'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').Value = '''0.07021'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').Value = '''1.005'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '''84.28'
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D18').Value = '''0.000009532'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').Value = '''16.71'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = '29.041.72'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').Value = '''5.339'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('D23').Value = '''10.98'
$ws.Range('E23').Value = '  +1.18%  '
$ws.Range('D24').Value = '''2.077'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '''157.42'
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = '''19.07'
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('D27').Value = '''5.649'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').Value = '''117.75'
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').Value = '''1.840'
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('E30').Value = '  +0.95%  '
$ws.Range('D31').Value = '''0.8593'
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('D32').Value = '''5.107'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = '''1.242'
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('D34').Value = '''3.020'
$ws.Range('E34').Value = '  +0.29%  '
$ws.Range('D35').Value = '''1.159'
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('D36').Value = '''0.05683'
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = '''1.003'
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.02045'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '''3.127'
$ws.Range('E39').Value = '  +15.80%  '
$ws.Range('D40').Value = '''7.438'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('D41').Value = '''0.5504'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').Value = '''0.1755'
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').Value = '''9.408'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').Value = '''2.199'
$ws.Range('E44').Value = '  +6.37%  '
$ws.Range('D45').Value = '''0.000002764'
$ws.Range('E45').Value = '  -7.38%  '
$ws.Range('D46').Value = '''0.5184'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').Value = '''11.20'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').Value = '''0.06908'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').Value = '''110.48'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('D51').Value = '''1.003'
$ws.Range('E51').Value = '  +0.24%  '
